$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.172.27"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.826.06"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'234.47"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "'0.6008"
$ws.Range("E6").Value = "  -3.87%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.06938"
$ws.Range("E8").Value = "  -5.84%  "
$ws.Range("D9").Value = "'0.2750"
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("D10").Value = "'23.41"
$ws.Range("E10").Value = "  -5.24%  "
$ws.Range("D12").Value = "1.831.72"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "'4.731"
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").Value = "'0.6233"
$ws.Range("E14").Value = "  -5.79%  "
$ws.Range("D15").Value = "'0.000009770"
$ws.Range("E15").Value = "  -7.87%  "
$ws.Range("D16").Value = "'77.37"
$ws.Range("E16").Value = "  -4.93%  "
$ws.Range("D17").Value = "28.877.45"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "'5.562"
$ws.Range("E18").Value = "  -10.71%  "
$ws.Range("D19").Value = "'216.38"
$ws.Range("E19").Value = "  -8.55%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'11.53"
$ws.Range("E21").Value = "  -5.55%  "
$ws.Range("D22").Value = "'6.882"
$ws.Range("E22").Value = "  -4.84%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'155.96"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "'7.920"
$ws.Range("E25").Value = "  -5.78%  "
$ws.Range("D26").Value = "'0.1288"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("D27").Value = "'16.46"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("D28").Value = "'0.06515"
$ws.Range("E28").Value = "  -6.84%  "
$ws.Range("D29").Value = "'1.409"
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("D30").Value = "'1.438"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").Value = "'3.822"
$ws.Range("E31").Value = "  -4.75%  "
$ws.Range("D32").Value = "'3.769"
$ws.Range("D33").Value = "'1.092"
$ws.Range("E33").Value = "  -5.16%  "
$ws.Range("D34").Value = "'1.722"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "'0.6448"
$ws.Range("E35").Value = "  -5.00%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'2.746"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'0.01755"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").Value = "'6.448"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").Value = "1.142.28"
$ws.Range("E40").Value = "  -7.21%  "
$ws.Range("D41").Value = "'0.8873"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "1.991.43"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "'100.23"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "'61.65"
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("D47").Value = "'1.607"
$ws.Range("E47").Value = "  -4.46%  "
$ws.Range("D48").Value = "'8.482"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").Value = "'0.05501"
$ws.Range("D50").Value = "'0.4538"
$ws.Range("D51").Value = "'6.390"
$ws.Range("E51").Value = "  -7.76%  "
